$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Data for the new Table1 on Sheet2 (Column1=RAND(), Column2=Suit, Column3=Number) ---
$rows = @(
    @("Spade", "Three"),
    @("Club", "Six"),
    @("Spade", "Five"),
    @("Club", "Three"),
    @("Spade", "Ten"),
    @("Diamond", "Six"),
    @("Heart", "Four"),
    @("Heart", "Ace"),
    @("Diamond", "Ace"),
    @("Diamond", "Four"),
    @("Diamond", "Three"),
    @("Spade", "Six"),
    @("Spade", "Queen"),
    @("Heart", "Ten"),
    @("Heart", "Three"),
    @("Spade", "Four"),
    @("Spade", "Two"),
    @("Spade", "Jack"),
    @("Spade", "Seven"),
    @("Club", "Seven"),
    @("Diamond", "Queen"),
    @("Club", "Queen"),
    @("Diamond", "Ten"),
    @("Spade", "King"),
    @("Spade", "Ace"),
    @("Heart", "Six"),
    @("Heart", "Two"),
    @("Club", "Two"),
    @("Club", "Four"),
    @("Diamond", "Eight"),
    @("Diamond", "King"),
    @("Diamond", "Jack"),
    @("Club", "Eight"),
    @("Heart", "Nine"),
    @("Club", "Nine"),
    @("Spade", "Nine"),
    @("Heart", "Five"),
    @("Spade", "Eight"),
    @("Club", "Ace"),
    @("Club", "Jack"),
    @("Diamond", "Nine"),
    @("Club", "King"),
    @("Heart", "Eight"),
    @("Diamond", "Two"),
    @("Heart", "King"),
    @("Club", "Ten"),
    @("Heart", "Queen"),
    @("Heart", "Jack"),
    @("Club", "Five"),
    @("Heart", "Seven"),
    @("Diamond", "Seven"),
    @("Diamond", "Five")
)

$ws2.Range("B2").Value = "Column1"
$ws2.Range("C2").Value = "Column2"
$ws2.Range("D2").Value = "Column3"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 3
    $ws2.Range("B$r").Formula = "=RAND()"
    $ws2.Range("C$r").Value = $rows[$i][0]
    $ws2.Range("D$r").Value = $rows[$i][1]
}

$ws2.Columns.Item(2).ColumnWidth = 10.166666666666666
$ws2.Columns.Item(3).ColumnWidth = 10.166666666666666
$ws2.Columns.Item(4).ColumnWidth = 10.166666666666666

# --- Turn B2:D54 into an Excel Table named Table1 ---
$tbl = $ws2.ListObjects.Add(1, $ws2.Range("B2:D54"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleLight9"

# --- Column F: CONCAT formula referencing the table ---
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 3
    $ws2.Range("F$r").Formula = '=_xlfn.CONCAT("cards.Add(new Card { Number=CardNumber.",Table1[[#This Row],[Column3]],", Suit=CardSuit.",Table1[[#This Row],[Column2]],"});")'
}

# --- G1 label cell ---
$ws2.Range("G1").Value = "cards.Add(new Card { Number = CardNumber.Ace, Suit = CardSuit.Heart });"

# --- Sheet view / selection swap: Sheet2 becomes the active/selected tab ---
$ws1.Range("A3:B54").Select()
$ws2.Range("F3:F54").Select()
